{"js": "// The \"tb_categoria\" table (2nd table in the document) loses its\n// \"classificacao\" and \"ativo\" attribute rows; those same two rows are\n// appended to the end of the \"tb_produto\" table (3rd table), with the\n// \"ativo\" row's justification text re-worded to talk about \"o produto\"\n// instead of \"a categoria\".\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst categoriaTable = tables.items[1];\nconst produtoTable = tables.items[2];\n\ncategoriaTable.rows.load(\"items\");\nawait context.sync();\n\n// Remove the last two rows (\"classificacao : varchar (255)\" and\n// \"ativo : boolean\") from the tb_categoria table.\nconst rows = categoriaTable.rows.items;\nrows[rows.length - 1].delete();\nrows[rows.length - 2].delete();\nawait context.sync();\n\n// Re-append the same two rows (attribute unchanged, reason text updated)\n// at the bottom of the tb_produto table.\nprodutoTable.addRows(Word.InsertLocation.end, 2, [\n  [\n    \"classificacao : varchar (255)\",\n    \"Classifica\u00e7\u00e3o do medicamento: Refer\u00eancia, Similar ou Gen\u00e9rico\",\n  ],\n  [\n    \"ativo : boolean\",\n    \"Atributo utilizado para ativar ou desativar o produto. Por exemplo: Desativar o produto x, pois o medicamento n\u00e3o ser\u00e1 distribu\u00eddo no posto de sa\u00fade por motivos quaisquer. \",\n  ],\n]);\nawait context.sync();\n", "ps1": "# The \"tb_categoria\" table (2nd table in the document) loses its\n# \"classificacao\" and \"ativo\" attribute rows; those same two rows are\n# appended to the end of the \"tb_produto\" table (3rd table), with the\n# \"ativo\" row's justification text re-worded to talk about \"o produto\"\n# instead of \"a categoria\".\n$d = $word.ActiveDocument\n\n$categoriaTable = $d.Tables.Item(2)\n\n# Remove the last two rows (\"classificacao : varchar (255)\" and\n# \"ativo : boolean\") from the tb_categoria table.\n$n = $categoriaTable.Rows.Count\n$categoriaTable.Rows.Item($n).Delete()\n$categoriaTable.Rows.Item($n - 1).Delete()\n\n# Re-resolve the tb_produto table AFTER the deletion above so its row\n# anchors reflect the now-shorter tb_categoria table.\n$produtoTable = $d.Tables.Item(3)\n\n# Re-append the same two rows (attribute unchanged, reason text updated)\n# at the bottom of the tb_produto table.\n$produtoTable.Rows.Add() | Out-Null\n$produtoTable.Rows.Add() | Out-Null\n$m = $produtoTable.Rows.Count\n\n$produtoTable.Cell($m - 1, 1).Range.Text = \"classificacao : varchar (255)\"\n$produtoTable.Cell($m - 1, 2).Range.Text = \"Classifica\u00e7\u00e3o do medicamento: Refer\u00eancia, Similar ou Gen\u00e9rico\"\n$produtoTable.Cell($m, 1).Range.Text = \"ativo : boolean\"\n$produtoTable.Cell($m, 2).Range.Text = \"Atributo utilizado para ativar ou desativar o produto. Por exemplo: Desativar o produto x, pois o medicamento n\u00e3o ser\u00e1 distribu\u00eddo no posto de sa\u00fade por motivos quaisquer. \"\n"}
